$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update price values in column D
$ws.Range("D33").Value = 652.812
$ws.Range("D34").Value = 1123.703
$ws.Range("D35").Value = 1242.108
